# Natmi following Dr Hou advice
# Replace the old 3-row (ECs/FAPs/sCs -> FAPs) result table with the new
# 8-row table (ECs/FAPs/Neutro/sCs sending clusters x ECs/FAPs target
# clusters) for the Agt -> Agtr2 ligand-receptor pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("ECs","Agt","Agtr2","ECs",2,1,0.257198,0.5143960000000001,0.1969864112529235,0.1476779584913764,1,0.5,0.1409635,0.281927,0.127302058387171,0.0886289202516707,0.03625553027300001,0.145022121092,0.02507677562679894,0.01308853800606173),
    @("ECs","Agt","Agtr2","FAPs",2,1,0.257198,0.5143960000000001,0.1969864112529235,0.1476779584913764,3,1,0.9663516666666667,2.899055,0.872697941612829,0.9113710797483293,0.2485437159633334,1.49126229578,0.1719096356261245,0.1345894204853147),
    @("FAPs","Agt","Agtr2","ECs",3,1,0.8072116666666668,2.421635,0.6182385918170351,0.6952272432352978,1,0.5,0.1409635,0.281927,0.127302058387171,0.0886289202516707,0.1137873817741667,0.6827242906450001,0.07870304531269458,0.06161723989749008),
    @("FAPs","Agt","Agtr2","FAPs",3,1,0.8072116666666668,2.421635,0.6182385918170351,0.6952272432352978,3,1,0.9663516666666667,2.899055,0.872697941612829,0.9113710797483293,0.7800503394361112,7.020453054925001,0.5395355465043405,0.6336100033378078),
    @("Neutro","Agt","Agtr2","ECs",1,0.3333333333333333,0.064689,0.194067,0.0495449185356821,0.05571469912391608,1,0.5,0.1409635,0.281927,0.127302058387171,0.0886289202516707,0.009118787851499999,0.05471272710899999,0.006307170112217033,0.004937933625499385),
    @("Neutro","Agt","Agtr2","FAPs",1,0.3333333333333333,0.064689,0.194067,0.0495449185356821,0.05571469912391608,3,1,0.9663516666666667,2.899055,0.872697941612829,0.9113710797483293,0.062512322965,0.562610906685,0.04323774842346507,0.05077676549841669),
    @("sCs","Agt","Agtr2","ECs",2,1,0.176565,0.35313,0.1352300783943593,0.1013800991494097,1,0.5,0.1409635,0.281927,0.127302058387171,0.0886289202516707,0.0248892203775,0.09955688151,0.01721506733546044,0.008985208722619497),
    @("sCs","Agt","Agtr2","FAPs",2,1,0.176565,0.35313,0.1352300783943593,0.1013800991494097,3,1,0.9663516666666667,2.899055,0.872697941612829,0.9113710797483293,0.170623882025,1.02374329215,0.1180150110588989,0.09239489042679017)
)

$r = 2
foreach ($rowdata in $rows) {
    $c = 1
    foreach ($val in $rowdata) {
        $ws.Cells.Item($r, $c).Value = $val
        $c = $c + 1
    }
    $r = $r + 1
}
